$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("AA2").Value = "2025-10-27 12:30:02"
$ws2.Range("AA3").Value = "2025-10-27 12:30:02"
$ws2.Range("AA4").Value = "2025-10-27 12:30:02"
$ws2.Range("AA5").Value = "2025-10-27 12:30:02"
$ws2.Range("AA6").Value = "2025-10-27 12:30:02"
$ws2.Range("AA7").Value = "2025-10-27 12:30:02"
$ws2.Range("AA8").Value = "2025-10-27 12:30:02"
$ws2.Range("C8").Value = 16
$ws2.Range("D8").Value = 312
$ws2.Range("E8").Value = 145
$ws2.Range("F8").Value = 167
$ws2.Range("G8").Value = 19.5
$ws2.Range("H8").Value = 9.06
$ws2.Range("I8").Value = 10.44
$ws2.Range("J8").Value = 65
$ws2.Range("K8").Value = 76
$ws2.Range("AA9").Value = "2025-10-27 12:30:02"
$ws2.Range("AA10").Value = "2025-10-27 12:30:02"
$ws2.Range("C10").Value = 13
$ws2.Range("D10").Value = 258
$ws2.Range("E10").Value = 126
$ws2.Range("F10").Value = 132
$ws2.Range("G10").Value = 19.85
$ws2.Range("H10").Value = 9.69
$ws2.Range("I10").Value = 10.15
$ws2.Range("J10").Value = 63
$ws2.Range("K10").Value = 56
$ws2.Range("W10").Value = 10
$ws2.Range("AA11").Value = "2025-10-27 12:30:02"
$ws2.Range("AA12").Value = "2025-10-27 12:30:02"
$ws2.Range("AA13").Value = "2025-10-27 12:30:02"
$ws2.Range("C13").Value = 7
$ws2.Range("D13").Value = 117
$ws2.Range("E13").Value = 60
$ws2.Range("F13").Value = 57
$ws2.Range("G13").Value = 16.71
$ws2.Range("H13").Value = 8.57
$ws2.Range("I13").Value = 8.140000000000001
$ws2.Range("J13").Value = 30
$ws2.Range("K13").Value = 26
$ws2.Range("AA14").Value = "2025-10-27 12:30:02"
$ws2.Range("AA15").Value = "2025-10-27 12:30:02"
$ws2.Range("C15").Value = 12
$ws2.Range("D15").Value = 178
$ws2.Range("E15").Value = 80
$ws2.Range("F15").Value = 98
$ws2.Range("G15").Value = 14.83
$ws2.Range("H15").Value = 6.67
$ws2.Range("I15").Value = 8.17
$ws2.Range("J15").Value = 40
$ws2.Range("K15").Value = 49
$ws2.Range("W15").Value = 8
$ws2.Range("AA16").Value = "2025-10-27 12:30:02"
$ws2.Range("AA17").Value = "2025-10-27 12:30:02"
$ws2.Range("AA18").Value = "2025-10-27 12:30:02"
$ws2.Range("C18").Value = 19
$ws2.Range("D18").Value = 243
$ws2.Range("E18").Value = 113
$ws2.Range("F18").Value = 130
$ws2.Range("G18").Value = 12.79
$ws2.Range("H18").Value = 5.95
$ws2.Range("I18").Value = 6.84
$ws2.Range("J18").Value = 54
$ws2.Range("K18").Value = 65
$ws2.Range("V18").Value = 6
$ws2.Range("AA19").Value = "2025-10-27 12:30:02"
$ws2.Range("C19").Value = 14
$ws2.Range("D19").Value = 261
$ws2.Range("E19").Value = 116
$ws2.Range("F19").Value = 145
$ws2.Range("G19").Value = 18.64
$ws2.Range("H19").Value = 8.289999999999999
$ws2.Range("I19").Value = 10.36
$ws2.Range("J19").Value = 53
$ws2.Range("K19").Value = 60
$ws2.Range("V19").Value = 6
$ws2.Range("AA20").Value = "2025-10-27 12:30:02"
$ws2.Range("AA21").Value = "2025-10-27 12:30:02"
$ws2.Range("AA22").Value = "2025-10-27 12:30:02"
$ws2.Range("AA23").Value = "2025-10-27 12:30:02"
$ws2.Range("AA24").Value = "2025-10-27 12:30:02"
$ws2.Range("AA25").Value = "2025-10-27 12:30:02"
$ws2.Range("AA26").Value = "2025-10-27 12:30:02"

$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("AA2").Value = "2025-10-27 12:30:02"
$ws3.Range("AA3").Value = "2025-10-27 12:30:02"
$ws3.Range("AA4").Value = "2025-10-27 12:30:02"
$ws3.Range("AA5").Value = "2025-10-27 12:30:02"
$ws3.Range("AA6").Value = "2025-10-27 12:30:02"
$ws3.Range("C6").Value = 11
$ws3.Range("D6").Value = 211
$ws3.Range("E6").Value = 88
$ws3.Range("F6").Value = 123
$ws3.Range("G6").Value = 19.18
$ws3.Range("H6").Value = 8
$ws3.Range("I6").Value = 11.18
$ws3.Range("J6").Value = 44
$ws3.Range("K6").Value = 59
$ws3.Range("V6").Value = 4
$ws3.Range("AA7").Value = "2025-10-27 12:30:02"
$ws3.Range("AA8").Value = "2025-10-27 12:30:02"
$ws3.Range("AA9").Value = "2025-10-27 12:30:02"
$ws3.Range("AA10").Value = "2025-10-27 12:30:02"
$ws3.Range("C10").Value = 11
$ws3.Range("D10").Value = 153
$ws3.Range("E10").Value = 72
$ws3.Range("F10").Value = 81
$ws3.Range("G10").Value = 13.91
$ws3.Range("H10").Value = 6.55
$ws3.Range("I10").Value = 7.36
$ws3.Range("J10").Value = 36
$ws3.Range("K10").Value = 38
$ws3.Range("W10").Value = 8
$ws3.Range("AA11").Value = "2025-10-27 12:30:02"
$ws3.Range("AA12").Value = "2025-10-27 12:30:02"
$ws3.Range("AA13").Value = "2025-10-27 12:30:02"
$ws3.Range("C13").Value = 15
$ws3.Range("D13").Value = 335
$ws3.Range("E13").Value = 169
$ws3.Range("F13").Value = 166
$ws3.Range("G13").Value = 22.33
$ws3.Range("H13").Value = 11.27
$ws3.Range("I13").Value = 11.07
$ws3.Range("J13").Value = 62
$ws3.Range("K13").Value = 63
$ws3.Range("AA14").Value = "2025-10-27 12:30:02"
$ws3.Range("AA15").Value = "2025-10-27 12:30:02"
$ws3.Range("AA16").Value = "2025-10-27 12:30:02"
$ws3.Range("AA17").Value = "2025-10-27 12:30:02"
$ws3.Range("C17").Value = 17
$ws3.Range("D17").Value = 317
$ws3.Range("E17").Value = 146
$ws3.Range("F17").Value = 171
$ws3.Range("G17").Value = 18.65
$ws3.Range("H17").Value = 8.59
$ws3.Range("I17").Value = 10.06
$ws3.Range("J17").Value = 68
$ws3.Range("K17").Value = 73
$ws3.Range("AA18").Value = "2025-10-27 12:30:02"
$ws3.Range("C18").Value = 13
$ws3.Range("D18").Value = 218
$ws3.Range("E18").Value = 111
$ws3.Range("F18").Value = 107
$ws3.Range("G18").Value = 16.77
$ws3.Range("H18").Value = 8.539999999999999
$ws3.Range("I18").Value = 8.23
$ws3.Range("J18").Value = 53
$ws3.Range("K18").Value = 51
$ws3.Range("AA19").Value = "2025-10-27 12:30:02"
$ws3.Range("C19").Value = 20
$ws3.Range("D19").Value = 436
$ws3.Range("E19").Value = 195
$ws3.Range("F19").Value = 241
$ws3.Range("G19").Value = 21.8
$ws3.Range("H19").Value = 9.75
$ws3.Range("I19").Value = 12.05
$ws3.Range("J19").Value = 90
$ws3.Range("K19").Value = 98
$ws3.Range("AA20").Value = "2025-10-27 12:30:02"
$ws3.Range("AA21").Value = "2025-10-27 12:30:02"
$ws3.Range("C21").Value = 11
$ws3.Range("D21").Value = 171
$ws3.Range("E21").Value = 80
$ws3.Range("F21").Value = 91
$ws3.Range("G21").Value = 15.55
$ws3.Range("H21").Value = 7.27
$ws3.Range("I21").Value = 8.27
$ws3.Range("J21").Value = 40
$ws3.Range("K21").Value = 43
$ws3.Range("V21").Value = 2
$ws3.Range("AA22").Value = "2025-10-27 12:30:02"
$ws3.Range("C22").Value = 19
$ws3.Range("D22").Value = 342
$ws3.Range("E22").Value = 147
$ws3.Range("F22").Value = 195
$ws3.Range("G22").Value = 18
$ws3.Range("H22").Value = 7.74
$ws3.Range("I22").Value = 10.26
$ws3.Range("J22").Value = 66
$ws3.Range("K22").Value = 85
$ws3.Range("W22").Value = 12
$ws3.Range("AA23").Value = "2025-10-27 12:30:02"
$ws3.Range("AA24").Value = "2025-10-27 12:30:02"
$ws3.Range("AA25").Value = "2025-10-27 12:30:02"
$ws3.Range("AA26").Value = "2025-10-27 12:30:02"
